$d = $word.ActiveDocument

$replacements = @(
    @("35×36=1260", "62×85=5270"),
    @("75×96=7200", "61×66=4026"),
    @("81×17=1377", "65×85=5525"),
    @("17×67=1139", "62×96=5952"),
    @("22×64=1408", "26×92=2392"),
    @("15×52=780",  "36×86=3096"),
    @("17×29=493",  "25×86=2150"),
    @("75×33=2475", "21×54=1134"),
    @("51×73=3723", "84×64=5376"),
    @("19×95=1805", "63×21=1323"),
    @("89×27=2403", "46×66=3036"),
    @("50×43=2150", "93×93=8649"),
    @("33×13=429",  "33×53=1749"),
    @("13×96=1248", "65×26=1690"),
    @("73×74=5402", "93×71=6603"),
    @("23×87=2001", "68×53=3604"),
    @("21×73=1533", "92×54=4968"),
    @("97×78=7566", "65×86=5590"),
    @("11×90=990",  "32×20=640"),
    @("52×79=4108", "47×27=1269"),
    @("22×48=1056", "76×20=1520"),
    @("30×72=2160", "43×47=2021"),
    @("56×48=2688", "22×55=1210"),
    @("54×80=4320", "78×54=4212"),
    @("94×15=1410", "82×89=7298")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
